# ST1 Line Setup - record new shift entries and normalize row 7 numeric columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns M..T (13..20) are numeric "product-characteristics4-7" columns in this sheet.
# Row 7 was left with those stored as text by a prior bad entry; fix it up to numeric,
# matching every other row in the log.
$numericCols = 13,14,15,16,17,18,19,20

foreach ($c in $numericCols) {
    $ws.Cells.Item(7, $c).Value = 4
}

# --- New shift-log rows appended below the existing data (rows 8-11) ---

function Set-TextCell($row, $col, $text) {
    # Force text storage even when the content looks numeric, matching the
    # "text" formatted columns (Date/Shift/VP/LT/OK-Remark/Sign) used throughout this log.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Set-NumericCell($row, $col, $number) {
    $ws.Cells.Item($row, $col).Value = $number
}

# Row 8
Set-TextCell 8 1 "2025-02-05T16:54"
Set-TextCell 8 2 "SHIFT_1"
Set-TextCell 8 3 "5644"
Set-TextCell 8 4 "654"
Set-TextCell 8 5 "65464"
Set-TextCell 8 6 "65464"
Set-TextCell 8 7 "OK"
Set-TextCell 8 8 "OK"
Set-TextCell 8 9 "OK"
Set-TextCell 8 10 "OK"
Set-TextCell 8 11 "OK"
Set-TextCell 8 12 "OK"
Set-NumericCell 8 13 654
Set-NumericCell 8 14 654
Set-NumericCell 8 15 64
Set-NumericCell 8 16 64
Set-NumericCell 8 17 64
Set-NumericCell 8 18 64
Set-NumericCell 8 19 646
Set-NumericCell 8 20 464
Set-TextCell 8 21 "OK"
Set-TextCell 8 22 "OK"
Set-TextCell 8 23 "654"
Set-TextCell 8 24 "654"
Set-TextCell 8 25 "654"

# Row 9
Set-TextCell 9 1 "2025-01-31T16:56"
Set-TextCell 9 2 "SHIFT_1"
Set-TextCell 9 3 "45"
Set-TextCell 9 4 "87"
Set-TextCell 9 5 "897"
Set-TextCell 9 6 "87"
Set-TextCell 9 7 "OK"
Set-TextCell 9 8 "OK"
Set-TextCell 9 9 "OK"
Set-TextCell 9 10 "OK"
Set-TextCell 9 11 "OK"
Set-TextCell 9 12 "OK"
Set-NumericCell 9 13 87
Set-NumericCell 9 14 87
Set-NumericCell 9 15 87
Set-NumericCell 9 16 87
Set-NumericCell 9 17 87
Set-NumericCell 9 18 87
Set-NumericCell 9 19 87
Set-NumericCell 9 20 87
Set-TextCell 9 21 "OK"
Set-TextCell 9 22 "OK"
Set-TextCell 9 23 "87"
Set-TextCell 9 24 "87"
Set-TextCell 9 25 "87"

# Row 10
Set-TextCell 10 1 "2025-02-04T17:12"
Set-TextCell 10 2 "SHIFT_1"
Set-TextCell 10 3 "5"
Set-TextCell 10 4 "5"
Set-TextCell 10 5 "5"
Set-TextCell 10 6 "5"
Set-TextCell 10 7 "OK"
Set-TextCell 10 8 "OK"
Set-TextCell 10 9 "OK"
Set-TextCell 10 10 "OK"
Set-TextCell 10 11 "OK"
Set-TextCell 10 12 "OK"
Set-NumericCell 10 13 5
Set-NumericCell 10 14 5
Set-NumericCell 10 15 5
Set-NumericCell 10 16 5
Set-NumericCell 10 17 5
Set-NumericCell 10 18 5
Set-NumericCell 10 19 5
Set-NumericCell 10 20 5
Set-TextCell 10 21 "OK"
Set-TextCell 10 22 "OK"
Set-TextCell 10 23 "5"
Set-TextCell 10 24 "5"
Set-TextCell 10 25 "5"

# Row 11 - note: unlike rows 8-10, columns M-T were NOT recorded numerically here
# (operator skipped the numeric-entry step), so they stay as plain text, same as
# row 7 did before it got its warning/fix-up.
Set-TextCell 11 1 "2025-02-06T17:27"
Set-TextCell 11 2 "SHIFT_1"
Set-TextCell 11 3 "45"
Set-TextCell 11 4 "45"
Set-TextCell 11 5 "54"
Set-TextCell 11 6 "45"
Set-TextCell 11 7 "OK"
Set-TextCell 11 8 "OK"
Set-TextCell 11 9 "OK"
Set-TextCell 11 10 "OK"
Set-TextCell 11 11 "OK"
Set-TextCell 11 12 "OK"
Set-TextCell 11 13 "974"
Set-TextCell 11 14 "8"
Set-TextCell 11 15 "45"
Set-TextCell 11 16 "45"
Set-TextCell 11 17 "45"
Set-TextCell 11 18 "54"
Set-TextCell 11 19 "54"
Set-TextCell 11 20 "45"
Set-TextCell 11 21 "OK"
Set-TextCell 11 22 "OK"
Set-TextCell 11 23 "54"
Set-TextCell 11 24 "54"
Set-TextCell 11 25 "54"
